# Insert a new data row at row 149 (pushing the existing row 149..216 down to 150..217),
# matching the semantics of the source diff: a new weekly price record was added and the
# remaining rows shifted down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new record's values.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are identical to the row that used to be
# at 149 (now shifted to 150); only D (Fecha) and the price/volume columns
# J, K, L, M, P carry new values per the diff.
$ws.Range("A149").Value = 11
$ws.Range("B149").Value = "Vega Monumental Concepción"
$ws.Range("C149").Value = "Bíobío"
$ws.Range("D149").Value = 44845
$ws.Range("E149").Value = 8
$ws.Range("F149").Value = 100112003
$ws.Range("G149").Value = "Ajo"
$ws.Range("H149").Value = "Chino"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 250
$ws.Range("K149").Value = 15000
$ws.Range("L149").Value = 16000
$ws.Range("M149").Value = 15600
$ws.Range("N149").Value = "$/caja 10 kilos"
$ws.Range("O149").Value = "China"
$ws.Range("P149").Value = 1560
$ws.Range("Q149").Value = 10
$ws.Range("R149").Value = "Hortaliza"

# Ensure the date cell keeps the date-formatted style (same as the other D-column cells,
# e.g. copy the number format from the row below which already carries style s="2").
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat
